# Update cryptos list values to reflect latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.356.04"
$ws.Range("E2").Value = "  -2.96%  "

# Row 3
$ws.Range("D3").Value = "2.463.82"
$ws.Range("E3").Value = "  -2.70%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.91"
$ws.Range("E5").Value = "  +0.74%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.08"
$ws.Range("E6").Value = "  -5.94%  "

# Row 7
$ws.Range("E7").Value = "  -3.31%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("E9").Value = "  -3.90%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.55"
$ws.Range("E10").Value = "  -6.59%  "

# Row 11
$ws.Range("E11").Value = "  -2.77%  "

# Row 12
$ws.Range("E12").Value = "  +0.04%  "

# Row 13
$ws.Range("E13").Value = "  -5.19%  "

# Row 14
$ws.Range("D14").Value = "2.839.92"
$ws.Range("E14").Value = "  -3.05%  "

# Row 15
$ws.Range("D15").Value = "2.479.26"
$ws.Range("E15").Value = "  -1.39%  "

# Row 16
$ws.Range("E16").Value = "  -8.67%  "

# Row 17
$ws.Range("E17").Value = "  -3.12%  "

# Row 18
$ws.Range("D18").Value = "41.340.74"
$ws.Range("E18").Value = "  -2.95%  "

# Row 19
$ws.Range("E19").Value = "  -6.01%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0917"
$ws.Range("E20").Value = "  -3.76%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.55"
$ws.Range("E21").Value = "  -5.02%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.90"
$ws.Range("E22").Value = "  -0.51%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.65"
$ws.Range("E23").Value = "  -2.18%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.77"
$ws.Range("E24").Value = "  -4.46%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.94"
$ws.Range("E25").Value = "  -4.92%  "

# Row 26
$ws.Range("E26").Value = "  +0.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.84"
$ws.Range("E27").Value = "  -4.10%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -5.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.73"
$ws.Range("E29").Value = "  -3.99%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.39"
$ws.Range("E30").Value = "  -6.94%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.41"
$ws.Range("E31").Value = "  -1.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.64"
$ws.Range("E32").Value = "  -2.10%  "

# Row 33
$ws.Range("E33").Value = "  -0.39%  "

# Row 34
$ws.Range("E34").Value = "  -7.61%  "

# Row 35
$ws.Range("E35").Value = "  -4.88%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.01"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.11"
$ws.Range("E37").Value = "  -7.40%  "

# Row 38
$ws.Range("E38").Value = "  -7.08%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.104"
$ws.Range("E39").Value = "  -5.89%  "

# Row 40
$ws.Range("E40").Value = "  -4.03%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.11"
$ws.Range("E41").Value = "  -5.24%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.32"
$ws.Range("E42").Value = "  -1.72%  "

# Row 43
$ws.Range("E43").Value = "  +0.10%  "

# Row 44
$ws.Range("D44").Value = "1.988.20"
$ws.Range("E44").Value = "  +1.00%  "

# Row 45
$ws.Range("E45").Value = "  -4.33%  "

# Row 46
$ws.Range("E46").Value = "  -7.58%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.77"
$ws.Range("E47").Value = "  -0.85%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "77.16"
$ws.Range("E48").Value = "  -4.87%  "

# Row 49
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.60"
$ws.Range("E49").Value = "  -3.41%  "

# Row 50
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.40"
$ws.Range("E50").Value = "  -4.08%  "

# Row 51
$ws.Range("E51").Value = "  -5.98%  "
